# Generate Report for Handoff
#
# The previous handoff (7dda93ed-dfce-488e-8073-c70f64c3b907.md) is now reported
# as "Ready for handoff" (instead of a failed transform) under a new id
# (fc1a2693-ec42-4cdf-a862-a22868795b5c.md), and a second file
# (ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md) was picked up and is also
# "Ready for handoff". The per-language sheets now show the generated
# handoff (.xlf) files together with their handoff datetime.

$wb = $excel.ActiveWorkbook

$mdFile1 = "fc1a2693-ec42-4cdf-a862-a22868795b5c.md"
$mdFile2 = "ffffef3041f4-79f5-4115-ba06-c94ca50d7b72.md"
$configFile = ".localization-config"

$mdFile1Url = "https://github.com/OpenLocalizationTest/oltest/blob/ba6004259b73dcc7b07181830d70be72dc883c71/e2e/" + $mdFile1
$mdFile2Url = "https://github.com/OpenLocalizationTest/oltest/blob/ba6004259b73dcc7b07181830d70be72dc883c71/e2e/" + $mdFile2
$configUrl = "https://github.com/OpenLocalizationTest/oltest/blob/bc8bea0acc59176eb4c6489517b7514d57e8f1c1/.localization-config"

$handoffBase = "fc1a2693-ec42-4cdf-a862-a22868795b5c.cdc48cec29ac5c32432d5e686fba457507182cc5"

# ---------------------------------------------------------------------------
# Overview sheet: File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Hyperlinks.Delete()

$overview.Cells.Item(2, 2).Value = "Ready for handoff"
$overview.Cells.Item(2, 3).Value = "Ready for handoff"

$overview.Cells.Item(3, 2).Value = "Ready for handoff"
$overview.Cells.Item(3, 3).Value = "Ready for handoff"

$overview.Cells.Item(4, 2).Value = "Not to be localized"
$overview.Cells.Item(4, 3).Value = "Not to be localized"

$overview.Hyperlinks.Add($overview.Range("A2"), $mdFile1Url, "", "", $mdFile1)
$overview.Hyperlinks.Add($overview.Range("A3"), $mdFile2Url, "", "", $mdFile2)
$overview.Hyperlinks.Add($overview.Range("A4"), $configUrl, "", "", $configFile)

# ---------------------------------------------------------------------------
# Per-language sheets (zh-cn, de-de)
# ---------------------------------------------------------------------------
$languages = @("zh-cn", "de-de")
$handoffDatetimes = @{ "zh-cn" = "2016-01-25 08:22:33"; "de-de" = "2016-01-25 08:22:45" }

foreach ($lang in $languages) {
  $ws = $wb.Worksheets.Item($lang)
  $ws.Hyperlinks.Delete()

  $handoffFile = $handoffBase + "." + $lang + ".xlf"
  $handoffFileUrl = "https://github.com/OpenLocalizationTest/oltest/blob/ba6004259b73dcc7b07181830d70be72dc883c71/e2e/handoff/" + $handoffFile
  $handoffDatetime = $handoffDatetimes[$lang]

  # Row 2: fc1a2693-...md, ready for handoff, handoff file + datetime
  $ws.Cells.Item(2, 2).Value = "Ready for handoff"
  $ws.Cells.Item(2, 4).Value = $handoffDatetime
  $ws.Cells.Item(2, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
  $ws.Cells.Item(2, 7).Value = "0001-01-01 00:00:00"
  $ws.Cells.Item(2, 8).Value = "Include"

  # Row 3: ffffef3041f4-...md, ready for handoff, same handoff file + datetime
  $ws.Cells.Item(3, 2).Value = "Ready for handoff"
  $ws.Cells.Item(3, 4).Value = $handoffDatetime
  $ws.Cells.Item(3, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
  $ws.Cells.Item(3, 7).Value = "0001-01-01 00:00:00"
  $ws.Cells.Item(3, 8).Value = "Include"

  # Row 4: .localization-config, not to be localized
  $ws.Cells.Item(4, 2).Value = "Not to be localized"
  $ws.Cells.Item(4, 4).Value = "0001-01-01 00:00:00"
  $ws.Cells.Item(4, 4).NumberFormat = "yyyy-mm-dd HH:mm:ss"
  $ws.Cells.Item(4, 7).Value = "0001-01-01 00:00:00"
  $ws.Cells.Item(4, 8).Value = "Ignored"

  $ws.Hyperlinks.Add($ws.Range("A2"), $mdFile1Url, "", "", $mdFile1)
  $ws.Hyperlinks.Add($ws.Range("C2"), $handoffFileUrl, "", "", $handoffFile)
  $ws.Hyperlinks.Add($ws.Range("A3"), $mdFile2Url, "", "", $mdFile2)
  $ws.Hyperlinks.Add($ws.Range("C3"), $handoffFileUrl, "", "", $handoffFile)
  $ws.Hyperlinks.Add($ws.Range("A4"), $configUrl, "", "", $configFile)
}
